$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Glyphs")

# Add the two new glyph rows (30 and 31) with their code / description pairs.
$ws.Range("A30").Value = "g29"
$ws.Range("B30").Value = "Superscript t"
$ws.Range("A31").Value = "g30"
$ws.Range("B31").Value = "q abbreviation"

# Make "Glyphs" the active sheet/tab and select the last entered cell, matching
# the new selection & active-tab state recorded in the workbook.
$ws.Activate()
$ws.Range("B31").Select()
